# Skip generation of volsmile and q-probability if already exists
#
# Updates the MarketObjects tracking sheet:
#   - 2025-09-05 (row 86) now also produced BTC.FUNDING.CSA_USD / BTCUSD.VOLSURFACE.REGULARIZED
#   - 2025-09-06 (row 87) now produced BTCUSD.SPOT
#   - appends rows for 2025-09-07 .. 2025-09-14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- amend existing rows -------------------------------------------------
$ws.Range("B86").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"
$ws.Range("B87").Value = "['BTCUSD.SPOT']"

# --- append new rows -------------------------------------------------------
$newRows = @(
    @("2025-09-07", "['BTCUSD.SPOT']"),
    @("2025-09-08", "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"),
    @("2025-09-09", "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"),
    @("2025-09-10", "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"),
    @("2025-09-11", "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"),
    @("2025-09-12", "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']"),
    @("2025-09-13", "[]"),
    @("2025-09-14", "[]")
)

$startRow = 88
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i

    # The date column holds plain text like "2025-09-07", not a real date
    # value. Force a text format before assigning so Excel's automatic
    # type-inference does not silently convert the string into a date
    # serial number, then drop the temporary format so the cell ends up
    # with the sheet's normal (unstyled) look.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newRows[$i][0]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
